$wb = $excel.ActiveWorkbook

# Add one additional day (2022-07-30) of violent-crime data across the
# Citywide Totals, By Neighborhood, and per-neighborhood detail sheets.

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 9).Value = 4041  # Aggravated Assault 2022: 4016 -> 4041
$ws.Cells.Item(3, 9).Value = 4203  # Aggravated Battery 2022: 4162 -> 4203
$ws.Cells.Item(4, 3).Value = 1807  # Criminal Sexual Assault 2016: 1808 -> 1807
$ws.Cells.Item(4, 8).Value = 1668  # Criminal Sexual Assault 2021: 1667 -> 1668
$ws.Cells.Item(4, 9).Value = 978  # Criminal Sexual Assault 2022: 974 -> 978
$ws.Cells.Item(5, 9).Value = 384  # Homicide 2022: 382 -> 384
$ws.Cells.Item(6, 9).Value = 4671  # Robbery 2022: 4645 -> 4671
$ws.Cells.Item(7, 3).Value = 28350  # Total 2016: 28351 -> 28350
$ws.Cells.Item(7, 8).Value = 25977  # Total 2021: 25976 -> 25977
$ws.Cells.Item(7, 9).Value = 14277  # Total 2022: 14179 -> 14277

# Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(4, 9).Value = 22  # Criminal Sexual Assault 2022: 21 -> 22
$ws.Cells.Item(7, 9).Value = 160  # Total 2022: 159 -> 160

# West Ridge
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(3, 9).Value = 44  # Aggravated Battery 2022: 43 -> 44
$ws.Cells.Item(7, 9).Value = 152  # Total 2022: 151 -> 152

# Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 9).Value = 148  # Aggravated Battery 2022: 146 -> 148
$ws.Cells.Item(7, 9).Value = 458  # Total 2022: 456 -> 458

# Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 9).Value = 73  # Aggravated Assault 2022: 72 -> 73
$ws.Cells.Item(3, 9).Value = 96  # Aggravated Battery 2022: 95 -> 96
$ws.Cells.Item(7, 9).Value = 271  # Total 2022: 269 -> 271

# North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 9).Value = 130  # Aggravated Assault 2022: 129 -> 130
$ws.Cells.Item(3, 9).Value = 199  # Aggravated Battery 2022: 198 -> 199
$ws.Cells.Item(4, 9).Value = 30  # Criminal Sexual Assault 2022: 29 -> 30
$ws.Cells.Item(7, 9).Value = 555  # Total 2022: 552 -> 555

# Gage Park
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(3, 9).Value = 37  # Aggravated Battery 2022: 36 -> 37
$ws.Cells.Item(7, 9).Value = 138  # Total 2022: 137 -> 138

# South Deering
$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(3, 9).Value = 41  # Aggravated Battery 2022: 40 -> 41
$ws.Cells.Item(7, 9).Value = 125  # Total 2022: 124 -> 125

# New City
$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(3, 9).Value = 92  # Aggravated Battery 2022: 90 -> 92
$ws.Cells.Item(6, 9).Value = 98  # Robbery 2022: 97 -> 98
$ws.Cells.Item(7, 9).Value = 322  # Total 2022: 319 -> 322

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(5, 9).Value = 45  # Armour Square 2022: 44 -> 45
$ws.Cells.Item(7, 9).Value = 452  # Auburn Gresham 2022: 448 -> 452
$ws.Cells.Item(8, 9).Value = 857  # Austin 2022: 852 -> 857
$ws.Cells.Item(9, 9).Value = 63  # Avalon Park 2022: 62 -> 63
$ws.Cells.Item(10, 9).Value = 95  # Avondale 2022: 94 -> 95
$ws.Cells.Item(11, 9).Value = 218  # Belmont Cragin 2022: 216 -> 218
$ws.Cells.Item(12, 9).Value = 31  # Beverly 2022: 28 -> 31
$ws.Cells.Item(15, 9).Value = 170  # Brighton Park 2022: 169 -> 170
$ws.Cells.Item(16, 9).Value = 39  # Bucktown 2022: 38 -> 39
$ws.Cells.Item(18, 9).Value = 100  # Calumet Heights 2022: 98 -> 100
$ws.Cells.Item(19, 9).Value = 396  # Chatham 2022: 392 -> 396
$ws.Cells.Item(20, 9).Value = 346  # Chicago Lawn 2022: 345 -> 346
$ws.Cells.Item(21, 9).Value = 79  # Chinatown 2022: 77 -> 79
$ws.Cells.Item(27, 9).Value = 131  # Edgewater 2022: 130 -> 131
$ws.Cells.Item(29, 9).Value = 915  # Englewood 2022: 912 -> 915
$ws.Cells.Item(31, 9).Value = 138  # Gage Park 2022: 137 -> 138
$ws.Cells.Item(33, 9).Value = 656  # Garfield Park 2022: 651 -> 656
$ws.Cells.Item(36, 9).Value = 197  # Grand Boulevard 2022: 196 -> 197
$ws.Cells.Item(37, 9).Value = 458  # Grand Crossing 2022: 456 -> 458
$ws.Cells.Item(41, 9).Value = 63  # Hermosa 2022: 62 -> 63
$ws.Cells.Item(42, 9).Value = 490  # Humboldt Park 2022: 489 -> 490
$ws.Cells.Item(44, 9).Value = 103  # Irving Park 2022: 102 -> 103
$ws.Cells.Item(48, 9).Value = 193  # Lake View 2022: 190 -> 193
$ws.Cells.Item(50, 9).Value = 61  # Lincoln Square 2022: 60 -> 61
$ws.Cells.Item(51, 9).Value = 142  # Little Italy, UIC 2022: 139 -> 142
$ws.Cells.Item(52, 9).Value = 308  # Little Village 2022: 306 -> 308
$ws.Cells.Item(53, 9).Value = 150  # Logan Square 2022: 149 -> 150
$ws.Cells.Item(54, 9).Value = 324  # Loop 2022: 322 -> 324
$ws.Cells.Item(55, 9).Value = 156  # Lower West Side 2022: 155 -> 156
$ws.Cells.Item(63, 3).Value = 240  # NO NEIGHBORHOOD DATA 2016: 241 -> 240
$ws.Cells.Item(63, 8).Value = 210  # NO NEIGHBORHOOD DATA 2021: 209 -> 210
$ws.Cells.Item(63, 9).Value = 56  # NO NEIGHBORHOOD DATA 2022: 53 -> 56
$ws.Cells.Item(64, 9).Value = 126  # Near South Side 2022: 125 -> 126
$ws.Cells.Item(65, 9).Value = 322  # New City 2022: 319 -> 322
$ws.Cells.Item(67, 9).Value = 555  # North Lawndale 2022: 552 -> 555
$ws.Cells.Item(68, 9).Value = 46  # North Park 2022: 45 -> 46
$ws.Cells.Item(72, 9).Value = 54  # Old Town 2022: 53 -> 54
$ws.Cells.Item(73, 9).Value = 117  # Portage Park 2022: 115 -> 117
$ws.Cells.Item(76, 9).Value = 213  # River North 2022: 212 -> 213
$ws.Cells.Item(77, 9).Value = 79  # Riverdale 2022: 78 -> 79
$ws.Cells.Item(79, 9).Value = 388  # Roseland 2022: 387 -> 388
$ws.Cells.Item(82, 9).Value = 18  # Sheffield & DePaul 2022: 17 -> 18
$ws.Cells.Item(83, 9).Value = 292  # South Chicago 2022: 290 -> 292
$ws.Cells.Item(84, 9).Value = 125  # South Deering 2022: 124 -> 125
$ws.Cells.Item(85, 9).Value = 642  # South Shore 2022: 639 -> 642
$ws.Cells.Item(88, 9).Value = 131  # United Center 2022: 128 -> 131
$ws.Cells.Item(89, 9).Value = 160  # Uptown 2022: 159 -> 160
$ws.Cells.Item(92, 9).Value = 43  # West Elsdon 2022: 42 -> 43
$ws.Cells.Item(93, 9).Value = 85  # West Lawn 2022: 84 -> 85
$ws.Cells.Item(94, 9).Value = 135  # West Loop 2022: 132 -> 135
$ws.Cells.Item(95, 9).Value = 234  # West Pullman 2022: 230 -> 234
$ws.Cells.Item(96, 9).Value = 152  # West Ridge 2022: 151 -> 152
$ws.Cells.Item(97, 9).Value = 106  # West Town 2022: 105 -> 106
$ws.Cells.Item(99, 9).Value = 271  # Woodlawn 2022: 269 -> 271
$ws.Cells.Item(101, 3).Value = 28350  # Total 2016: 28351 -> 28350
$ws.Cells.Item(101, 8).Value = 25977  # Total 2021: 25976 -> 25977
$ws.Cells.Item(101, 9).Value = 14277  # Total 2022: 14179 -> 14277

# South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(3, 9).Value = 113  # Aggravated Battery 2022: 112 -> 113
$ws.Cells.Item(4, 9).Value = 13  # Criminal Sexual Assault 2022: 12 -> 13
$ws.Cells.Item(7, 9).Value = 292  # Total 2022: 290 -> 292

# West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 9).Value = 84  # Aggravated Assault 2022: 82 -> 84
$ws.Cells.Item(3, 9).Value = 89  # Aggravated Battery 2022: 88 -> 89
$ws.Cells.Item(6, 9).Value = 39  # Robbery 2022: 38 -> 39
$ws.Cells.Item(7, 9).Value = 234  # Total 2022: 230 -> 234

# Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 9).Value = 244  # Aggravated Battery 2022: 241 -> 244
$ws.Cells.Item(6, 9).Value = 202  # Robbery 2022: 200 -> 202
$ws.Cells.Item(7, 9).Value = 656  # Total 2022: 651 -> 656

# Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 9).Value = 75  # Aggravated Assault 2022: 74 -> 75
$ws.Cells.Item(6, 9).Value = 163  # Robbery 2022: 162 -> 163
$ws.Cells.Item(7, 9).Value = 324  # Total 2022: 322 -> 324

# Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(3, 9).Value = 315  # Aggravated Battery 2022: 313 -> 315
$ws.Cells.Item(6, 9).Value = 250  # Robbery 2022: 249 -> 250
$ws.Cells.Item(7, 9).Value = 915  # Total 2022: 912 -> 915

# Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 9).Value = 145  # Aggravated Assault 2022: 144 -> 145
$ws.Cells.Item(3, 9).Value = 116  # Aggravated Battery 2022: 113 -> 116
$ws.Cells.Item(7, 9).Value = 396  # Total 2022: 392 -> 396

# Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(2, 9).Value = 35  # Aggravated Assault 2022: 34 -> 35
$ws.Cells.Item(7, 9).Value = 103  # Total 2022: 102 -> 103

# Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(6, 9).Value = 114  # Robbery 2022: 111 -> 114
$ws.Cells.Item(7, 9).Value = 193  # Total 2022: 190 -> 193

# River North
$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(3, 9).Value = 51  # Aggravated Battery 2022: 50 -> 51
$ws.Cells.Item(7, 9).Value = 213  # Total 2022: 212 -> 213

# South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 9).Value = 167  # Aggravated Assault 2022: 165 -> 167
$ws.Cells.Item(6, 9).Value = 158  # Robbery 2022: 157 -> 158
$ws.Cells.Item(7, 9).Value = 642  # Total 2022: 639 -> 642

# Hermosa
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Cells.Item(3, 9).Value = 22  # Aggravated Battery 2022: 21 -> 22
$ws.Cells.Item(7, 9).Value = 63  # Total 2022: 62 -> 63

# Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(3, 9).Value = 168  # Aggravated Battery 2022: 167 -> 168
$ws.Cells.Item(7, 9).Value = 490  # Total 2022: 489 -> 490

# Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(2, 9).Value = 33  # Aggravated Assault 2022: 32 -> 33
$ws.Cells.Item(7, 9).Value = 95  # Total 2022: 94 -> 95

# Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(3, 9).Value = 45  # Aggravated Battery 2022: 44 -> 45
$ws.Cells.Item(7, 9).Value = 156  # Total 2022: 155 -> 156

# Chinatown
$ws = $wb.Worksheets.Item('Chinatown')
$ws.Cells.Item(3, 9).Value = 11  # Aggravated Battery 2022: 9 -> 11
$ws.Cells.Item(7, 9).Value = 79  # Total 2022: 77 -> 79

# Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 9).Value = 122  # Aggravated Battery 2022: 121 -> 122
$ws.Cells.Item(7, 9).Value = 388  # Total 2022: 387 -> 388

# Near South Side
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(2, 9).Value = 34  # Aggravated Assault 2022: 33 -> 34
$ws.Cells.Item(7, 9).Value = 126  # Total 2022: 125 -> 126

# Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(3, 9).Value = 106  # Aggravated Battery 2022: 105 -> 106
$ws.Cells.Item(7, 9).Value = 346  # Total 2022: 345 -> 346

# Calumet Heights
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(2, 9).Value = 30  # Aggravated Assault 2022: 29 -> 30
$ws.Cells.Item(3, 9).Value = 24  # Aggravated Battery 2022: 23 -> 24
$ws.Cells.Item(7, 9).Value = 100  # Total 2022: 98 -> 100

# Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(2, 9).Value = 63  # Aggravated Assault 2022: 62 -> 63
$ws.Cells.Item(7, 9).Value = 197  # Total 2022: 196 -> 197

# West Lawn
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(6, 9).Value = 37  # Robbery 2022: 36 -> 37
$ws.Cells.Item(7, 9).Value = 85  # Total 2022: 84 -> 85

# Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 9).Value = 88  # Aggravated Assault 2022: 86 -> 88
$ws.Cells.Item(7, 9).Value = 308  # Total 2022: 306 -> 308

# West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(6, 9).Value = 76  # Robbery 2022: 73 -> 76
$ws.Cells.Item(7, 9).Value = 135  # Total 2022: 132 -> 135

# Brighton Park
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(3, 9).Value = 38  # Aggravated Battery 2022: 37 -> 38
$ws.Cells.Item(7, 9).Value = 170  # Total 2022: 169 -> 170

# Lincoln Square
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(6, 9).Value = 19  # Robbery 2022: 18 -> 19
$ws.Cells.Item(7, 9).Value = 61  # Total 2022: 60 -> 61

# Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(6, 9).Value = 55  # Robbery 2022: 53 -> 55
$ws.Cells.Item(7, 9).Value = 218  # Total 2022: 216 -> 218

# Avalon Park
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(2, 9).Value = 25  # Aggravated Assault 2022: 24 -> 25
$ws.Cells.Item(7, 9).Value = 63  # Total 2022: 62 -> 63

# Portage Park
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 9).Value = 43  # Aggravated Assault 2022: 41 -> 43
$ws.Cells.Item(7, 9).Value = 117  # Total 2022: 115 -> 117

# West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(3, 9).Value = 20  # Aggravated Battery 2022: 19 -> 20
$ws.Cells.Item(7, 9).Value = 106  # Total 2022: 105 -> 106

# West Elsdon
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Cells.Item(6, 9).Value = 18  # Robbery 2022: 17 -> 18
$ws.Cells.Item(7, 9).Value = 43  # Total 2022: 42 -> 43

# United Center
$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 9).Value = 35  # Aggravated Assault 2022: 33 -> 35
$ws.Cells.Item(3, 9).Value = 46  # Aggravated Battery 2022: 45 -> 46
$ws.Cells.Item(7, 9).Value = 131  # Total 2022: 128 -> 131

# Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 9).Value = 267  # Aggravated Assault 2022: 266 -> 267
$ws.Cells.Item(3, 9).Value = 241  # Aggravated Battery 2022: 237 -> 241
$ws.Cells.Item(7, 9).Value = 857  # Total 2022: 852 -> 857

# Armour Square
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Cells.Item(6, 9).Value = 22  # Robbery 2022: 21 -> 22
$ws.Cells.Item(7, 9).Value = 45  # Total 2022: 44 -> 45

# Edgewater
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(6, 9).Value = 52  # Robbery 2022: 51 -> 52
$ws.Cells.Item(7, 9).Value = 131  # Total 2022: 130 -> 131

# Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(3, 9).Value = 43  # Aggravated Battery 2022: 42 -> 43
$ws.Cells.Item(4, 9).Value = 16  # Criminal Sexual Assault 2022: 15 -> 16
$ws.Cells.Item(6, 9).Value = 54  # Robbery 2022: 53 -> 54
$ws.Cells.Item(7, 9).Value = 142  # Total 2022: 139 -> 142

# North Park
$ws = $wb.Worksheets.Item('North Park')
$ws.Cells.Item(4, 9).Value = 7  # Criminal Sexual Assault 2022: 6 -> 7
$ws.Cells.Item(7, 9).Value = 46  # Total 2022: 45 -> 46

# Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(2, 9).Value = 31  # Aggravated Assault 2022: 30 -> 31
$ws.Cells.Item(7, 9).Value = 150  # Total 2022: 149 -> 150

# Old Town
$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(6, 9).Value = 31  # Robbery 2022: 30 -> 31
$ws.Cells.Item(7, 9).Value = 54  # Total 2022: 53 -> 54

# Sheffield & DePaul
$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Cells.Item(5, 9).Value = 11  # Robbery 2022: 10 -> 11
$ws.Cells.Item(6, 9).Value = 18  # Total 2022: 17 -> 18

# Riverdale
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Cells.Item(2, 9).Value = 25  # Aggravated Assault 2022: 24 -> 25
$ws.Cells.Item(7, 9).Value = 79  # Total 2022: 78 -> 79

# Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 9).Value = 155  # Aggravated Assault 2022: 153 -> 155
$ws.Cells.Item(3, 9).Value = 141  # Aggravated Battery 2022: 139 -> 141
$ws.Cells.Item(7, 9).Value = 452  # Total 2022: 448 -> 452

# Beverly
$ws = $wb.Worksheets.Item('Beverly')
$ws.Cells.Item(5, 9).Value = 1  # Homicide 2022: (blank) -> 1
$ws.Cells.Item(6, 9).Value = 14  # Robbery 2022: 12 -> 14
$ws.Cells.Item(7, 9).Value = 31  # Total 2022: 28 -> 31

# Bucktown
$ws = $wb.Worksheets.Item('Bucktown')
$ws.Cells.Item(5, 9).Value = 1  # Homicide 2022: (blank) -> 1
$ws.Cells.Item(7, 9).Value = 39  # Total 2022: 38 -> 39
